$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 121; this shifts the existing rows 121-229
# down to 122-230 and extends the used range to A1:T230.
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new weekly record.
$ws.Range("A121").Value = 4
$ws.Range("B121").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C121").Value = "Los Lagos"
$ws.Range("D121").Value = 44789
$ws.Range("E121").Value = 10
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100108
$ws.Range("H121").Value = "Tropicales y subtropicales"
$ws.Range("I121").Value = 100108002
$ws.Range("J121").Value = "Mango"
$ws.Range("K121").Value = "Sin especificar"
$ws.Range("L121").Value = "Primera"
$ws.Range("M121").Value = 160
$ws.Range("N121").Value = 14000
$ws.Range("O121").Value = 14000
$ws.Range("P121").Value = 14000
$ws.Range("Q121").Value = "`$/bandeja 4 kilos"
$ws.Range("R121").Value = "Brasil"
$ws.Range("S121").Value = 3500
$ws.Range("T121").Value = 4
